# Existing ICDC Biobank filter fixes
# Fix the SamplesTab query (cell B3 on the "startup" sheet):
#  - trim stray trailing spaces on several lines
#  - wrap samp.specific_sample_pathology in replace(..., "  ", " ") to
#    collapse double spaces in the Pathology/Morphology column

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis)
MATCH (r:registration)-->(c)
WHERE s.clinical_study_designation IN ['MGT01'] and r.registration_origin in['CSU ACTR','UCD SVM']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS ``Sample ID``,
        coalesce(c.case_id, '') AS ``Case ID``,
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis,
        coalesce(samp.sample_site, '') AS ``Sample Site``,
        coalesce(samp.summarized_sample_type, '') AS ``Sample Type``,
        replace(coalesce(samp.specific_sample_pathology, ''), `"  `", `" `") AS ``Pathology/Morphology``,
        coalesce(samp.tumor_grade, '') AS ``Tumor Grade``,
        coalesce(samp.sample_chronology, '') AS ``Sample Chronology``,
        coalesce(samp.percentage_tumor, '') AS ``Percentage Tumor``,
        coalesce(samp.necropsy_sample, '') AS ``Necropsy Sample``,
        coalesce(samp.sample_preservation, '') AS ``Sample Preservation``
order by samp.sample_id asc
limit 100
"@

# Here-strings keep a trailing newline before the closing "@ marker; strip it
# so the cell text ends exactly at "limit 100" with no trailing newline.
$newQuery = $newQuery.TrimEnd("`r", "`n")

$ws.Range("B3").Value = $newQuery
